$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update footer "last updated" timestamp text (A1)
$ws.Range("A1").Value = "Datos actualizados a 4 de Mayo de 2020 a las 14:33"

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 1188898
$ws.Cells.Item(4, 3).Value = 776
$ws.Cells.Item(4, 5).Value = 941696
$ws.Cells.Item(4, 7).Value = 10
$ws.Cells.Item(4, 8).Value = 68608

# Row 18: India
$ws.Cells.Item(18, 2).Value = 42836
$ws.Cells.Item(18, 3).Value = 331
$ws.Cells.Item(18, 5).Value = 29659

# Row 19: Paises Bajos
$ws.Cells.Item(19, 2).Value = 40770
$ws.Cells.Item(19, 3).Value = 199
$ws.Cells.Item(19, 5).Value = 35438
$ws.Cells.Item(19, 7).Value = 26
$ws.Cells.Item(19, 8).Value = 5082

# Row 42: Dinamarca
$ws.Cells.Item(42, 4).Value = 7088
$ws.Cells.Item(42, 5).Value = 2089
$ws.Cells.Item(42, 6).Value = 57
$ws.Cells.Item(42, 7).Value = 9
$ws.Cells.Item(42, 8).Value = 493

# Row 46: Noruega
$ws.Cells.Item(46, 6).Value = 27

# Row 65: Afganistan
$ws.Cells.Item(65, 2).Value = 2894
$ws.Cells.Item(65, 3).Value = 190
$ws.Cells.Item(65, 4).Value = 397
$ws.Cells.Item(65, 5).Value = 2407
$ws.Cells.Item(65, 7).Value = 5
$ws.Cells.Item(65, 8).Value = 90

# Rows 71-72: Uzbekistan and Ghana swap order, with Uzbekistan's stats updated
$ws.Cells.Item(71, 1).Value = "Uzbekistan"
$ws.Cells.Item(71, 2).Value = 2181
$ws.Cells.Item(71, 3).Value = 32
$ws.Cells.Item(71, 4).Value = 1370
$ws.Cells.Item(71, 5).Value = 801
$ws.Cells.Item(71, 6).Value = 8
$ws.Cells.Item(71, 7).Value = 0
$ws.Cells.Item(71, 8).Value = 10

$ws.Cells.Item(72, 1).Value = "Ghana"
$ws.Cells.Item(72, 2).Value = 2169
$ws.Cells.Item(72, 3).Value = 0
$ws.Cells.Item(72, 4).Value = 229
$ws.Cells.Item(72, 5).Value = 1922
$ws.Cells.Item(72, 6).Value = 4
$ws.Cells.Item(72, 7).Value = 0
$ws.Cells.Item(72, 8).Value = 18

# Row 73: Croacia stats updated
$ws.Cells.Item(73, 2).Value = 2101
$ws.Cells.Item(73, 3).Value = 5
$ws.Cells.Item(73, 4).Value = 1522
$ws.Cells.Item(73, 5).Value = 499
$ws.Cells.Item(73, 7).Value = 1
$ws.Cells.Item(73, 8).Value = 80
